$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 5's formatting (incl. the date-style cell G) down into row 6
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial()

# Now set the actual values for the new row
$ws.Range("A6").Value = 10051.719999999999
$ws.Range("B6").Value = 9928.61
$ws.Range("C6").Value = 19.36
$ws.Range("D6").Value = 19.12
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = -1.24
$ws.Range("G6").Value = 42612.674780092595
$ws.Range("H6").Value = $true
